# Append a new attendance row (row 3) to the Attendance sheet:
# Rishabh | 2025-12-19 | Friday | 18:23:59
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Rishabh"

# B3 looks like a date ("2025-12-19"); Excel would normally auto-convert it
# to a date serial number. Force it to stay a plain text string (matching
# the other text cells in the sheet) by marking the cell as Text before
# assigning the value...
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-12-19"
# ...then re-apply the plain/default style (same as the neighbouring cell)
# so no stray number-format style lingers on the cell.
$ws.Range("B3").Style = $ws.Range("A3").Style

$ws.Range("C3").Value = "Friday"
$ws.Range("D3").Value = "18:23:59"
